$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) in the s_vals sheet, reusing the same header
# formatting (bold, bordered, centered) already applied to the other
# header cells (e.g. G1) by copying its format onto H1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the corresponding data row with the computed "Save" value.
$ws.Range("H2").Value = 0
